$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 0.7922663333333334
$ws.Cells.Item(2, 8).Value = 2.376799
$ws.Cells.Item(2, 9).Value = 0.1759587713796512
$ws.Cells.Item(2, 10).Value = 0.1759587713796512
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 24.365583
$ws.Cells.Item(2, 14).Value = 73.096749
$ws.Cells.Item(2, 15).Value = 0.3097154004536173
$ws.Cells.Item(2, 16).Value = 0.3097154004536173
$ws.Cells.Item(2, 17).Value = 19.304031102939
$ws.Cells.Item(2, 18).Value = 173.736279926451
$ws.Cells.Item(2, 19).Value = 0.05449714134117516
$ws.Cells.Item(2, 20).Value = 0.05449714134117517

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 0.7922663333333334
$ws.Cells.Item(3, 8).Value = 2.376799
$ws.Cells.Item(3, 9).Value = 0.1759587713796512
$ws.Cells.Item(3, 10).Value = 0.1759587713796512
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.228158333333333
$ws.Cells.Item(3, 14).Value = 3.684475
$ws.Cells.Item(3, 15).Value = 0.0156113461364245
$ws.Cells.Item(3, 16).Value = 0.0156113461364245
$ws.Cells.Item(3, 17).Value = 0.9730284995027778
$ws.Cells.Item(3, 18).Value = 8.757256495525001
$ws.Cells.Item(3, 19).Value = 0.00274695328574772
$ws.Cells.Item(3, 20).Value = 0.00274695328574772

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 0.7922663333333334
$ws.Cells.Item(4, 8).Value = 2.376799
$ws.Cells.Item(4, 9).Value = 0.1759587713796512
$ws.Cells.Item(4, 10).Value = 0.1759587713796512
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 48.48145033333333
$ws.Cells.Item(4, 14).Value = 145.444351
$ws.Cells.Item(4, 15).Value = 0.6162566192058893
$ws.Cells.Item(4, 16).Value = 0.6162566192058893
$ws.Cells.Item(4, 17).Value = 38.41022089027211
$ws.Cells.Item(4, 18).Value = 345.691988012449
$ws.Cells.Item(4, 19).Value = 0.1084357575700458
$ws.Cells.Item(4, 20).Value = 0.1084357575700458

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.7922663333333334
$ws.Cells.Item(5, 8).Value = 2.376799
$ws.Cells.Item(5, 9).Value = 0.1759587713796512
$ws.Cells.Item(5, 10).Value = 0.1759587713796512
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 4.595688
$ws.Cells.Item(5, 14).Value = 13.787064
$ws.Cells.Item(5, 15).Value = 0.05841663420406906
$ws.Cells.Item(5, 16).Value = 0.05841663420406905
$ws.Cells.Item(5, 17).Value = 3.641008880904
$ws.Cells.Item(5, 18).Value = 32.769079928136
$ws.Cells.Item(5, 19).Value = 0.0102789191826825
$ws.Cells.Item(5, 20).Value = 0.0102789191826825

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 2.077831666666667
$ws.Cells.Item(6, 8).Value = 6.233495
$ws.Cells.Item(6, 9).Value = 0.4614770208171574
$ws.Cells.Item(6, 10).Value = 0.4614770208171574
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 24.365583
$ws.Cells.Item(6, 14).Value = 73.096749
$ws.Cells.Item(6, 15).Value = 0.3097154004536173
$ws.Cells.Item(6, 16).Value = 0.3097154004536173
$ws.Cells.Item(6, 17).Value = 50.62757993419501
$ws.Cells.Item(6, 18).Value = 455.6482194077551
$ws.Cells.Item(6, 19).Value = 0.1429265403025282
$ws.Cells.Item(6, 20).Value = 0.1429265403025282

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 2.077831666666667
$ws.Cells.Item(7, 8).Value = 6.233495
$ws.Cells.Item(7, 9).Value = 0.4614770208171574
$ws.Cells.Item(7, 10).Value = 0.4614770208171574
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.228158333333333
$ws.Cells.Item(7, 14).Value = 3.684475
$ws.Cells.Item(7, 15).Value = 0.0156113461364245
$ws.Cells.Item(7, 16).Value = 0.0156113461364245
$ws.Cells.Item(7, 17).Value = 2.551906276680556
$ws.Cells.Item(7, 18).Value = 22.967156490125
$ws.Cells.Item(7, 19).Value = 0.007204277505982619
$ws.Cells.Item(7, 20).Value = 0.007204277505982619

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.077831666666667
$ws.Cells.Item(8, 8).Value = 6.233495
$ws.Cells.Item(8, 9).Value = 0.4614770208171574
$ws.Cells.Item(8, 10).Value = 0.4614770208171574
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 48.48145033333333
$ws.Cells.Item(8, 14).Value = 145.444351
$ws.Cells.Item(8, 15).Value = 0.6162566192058893
$ws.Cells.Item(8, 16).Value = 0.6162566192058893
$ws.Cells.Item(8, 17).Value = 100.7362927485272
$ws.Cells.Item(8, 18).Value = 906.626634736745
$ws.Cells.Item(8, 19).Value = 0.2843882686899872
$ws.Cells.Item(8, 20).Value = 0.2843882686899872

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.077831666666667
$ws.Cells.Item(9, 8).Value = 6.233495
$ws.Cells.Item(9, 9).Value = 0.4614770208171574
$ws.Cells.Item(9, 10).Value = 0.4614770208171574
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.595688
$ws.Cells.Item(9, 14).Value = 13.787064
$ws.Cells.Item(9, 15).Value = 0.05841663420406906
$ws.Cells.Item(9, 16).Value = 0.05841663420406905
$ws.Cells.Item(9, 17).Value = 9.549066056520001
$ws.Cells.Item(9, 18).Value = 85.94159450868
$ws.Cells.Item(9, 19).Value = 0.02695793431865944
$ws.Cells.Item(9, 20).Value = 0.02695793431865944

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.8480786666666668
$ws.Cells.Item(10, 8).Value = 2.544236
$ws.Cells.Item(10, 9).Value = 0.1883544383264543
$ws.Cells.Item(10, 10).Value = 0.1883544383264543
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 24.365583
$ws.Cells.Item(10, 14).Value = 73.096749
$ws.Cells.Item(10, 15).Value = 0.3097154004536173
$ws.Cells.Item(10, 16).Value = 0.3097154004536173
$ws.Cells.Item(10, 17).Value = 20.663931143196
$ws.Cells.Item(10, 18).Value = 185.975380288764
$ws.Cells.Item(10, 19).Value = 0.05833627029349394
$ws.Cells.Item(10, 20).Value = 0.05833627029349395

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.8480786666666668
$ws.Cells.Item(11, 8).Value = 2.544236
$ws.Cells.Item(11, 9).Value = 0.1883544383264543
$ws.Cells.Item(11, 10).Value = 0.1883544383264543
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 1.228158333333333
$ws.Cells.Item(11, 14).Value = 3.684475
$ws.Cells.Item(11, 15).Value = 0.0156113461364245
$ws.Cells.Item(11, 16).Value = 0.0156113461364245
$ws.Cells.Item(11, 17).Value = 1.041574881788889
$ws.Cells.Item(11, 18).Value = 9.3741739361
$ws.Cells.Item(11, 19).Value = 0.002940466333046099
$ws.Cells.Item(11, 20).Value = 0.002940466333046099

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.8480786666666668
$ws.Cells.Item(12, 8).Value = 2.544236
$ws.Cells.Item(12, 9).Value = 0.1883544383264543
$ws.Cells.Item(12, 10).Value = 0.1883544383264543
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 48.48145033333333
$ws.Cells.Item(12, 14).Value = 145.444351
$ws.Cells.Item(12, 15).Value = 0.6162566192058893
$ws.Cells.Item(12, 16).Value = 0.6162566192058893
$ws.Cells.Item(12, 17).Value = 41.11608375675956
$ws.Cells.Item(12, 18).Value = 370.044753810836
$ws.Cells.Item(12, 19).Value = 0.1160746693754849
$ws.Cells.Item(12, 20).Value = 0.1160746693754849

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.8480786666666668
$ws.Cells.Item(13, 8).Value = 2.544236
$ws.Cells.Item(13, 9).Value = 0.1883544383264543
$ws.Cells.Item(13, 10).Value = 0.1883544383264543
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 4.595688
$ws.Cells.Item(13, 14).Value = 13.787064
$ws.Cells.Item(13, 15).Value = 0.05841663420406906
$ws.Cells.Item(13, 16).Value = 0.05841663420406905
$ws.Cells.Item(13, 17).Value = 3.897504951456
$ws.Cells.Item(13, 18).Value = 35.077544563104
$ws.Cells.Item(13, 19).Value = 0.01100303232442936
$ws.Cells.Item(13, 20).Value = 0.01100303232442936

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 0.7843913333333336
$ws.Cells.Item(14, 8).Value = 2.353174000000001
$ws.Cells.Item(14, 9).Value = 0.1742097694767371
$ws.Cells.Item(14, 10).Value = 0.1742097694767372
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 24.365583
$ws.Cells.Item(14, 14).Value = 73.096749
$ws.Cells.Item(14, 15).Value = 0.3097154004536173
$ws.Cells.Item(14, 16).Value = 0.3097154004536173
$ws.Cells.Item(14, 17).Value = 19.11215213681401
$ws.Cells.Item(14, 18).Value = 172.009369231326
$ws.Cells.Item(14, 19).Value = 0.05395544851641999
$ws.Cells.Item(14, 20).Value = 0.05395544851642001

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 0.7843913333333336
$ws.Cells.Item(15, 8).Value = 2.353174000000001
$ws.Cells.Item(15, 9).Value = 0.1742097694767371
$ws.Cells.Item(15, 10).Value = 0.1742097694767372
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 1.228158333333333
$ws.Cells.Item(15, 14).Value = 3.684475
$ws.Cells.Item(15, 15).Value = 0.0156113461364245
$ws.Cells.Item(15, 16).Value = 0.0156113461364245
$ws.Cells.Item(15, 17).Value = 0.963356752627778
$ws.Cells.Item(15, 18).Value = 8.670210773650002
$ws.Cells.Item(15, 19).Value = 0.002719649011648063
$ws.Cells.Item(15, 20).Value = 0.002719649011648063

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 0.7843913333333336
$ws.Cells.Item(16, 8).Value = 2.353174000000001
$ws.Cells.Item(16, 9).Value = 0.1742097694767371
$ws.Cells.Item(16, 10).Value = 0.1742097694767372
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 48.48145033333333
$ws.Cells.Item(16, 14).Value = 145.444351
$ws.Cells.Item(16, 15).Value = 0.6162566192058893
$ws.Cells.Item(16, 16).Value = 0.6162566192058893
$ws.Cells.Item(16, 17).Value = 38.02842946889712
$ws.Cells.Item(16, 18).Value = 342.255865220074
$ws.Cells.Item(16, 19).Value = 0.1073579235703713
$ws.Cells.Item(16, 20).Value = 0.1073579235703714

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 0.7843913333333336
$ws.Cells.Item(17, 8).Value = 2.353174000000001
$ws.Cells.Item(17, 9).Value = 0.1742097694767371
$ws.Cells.Item(17, 10).Value = 0.1742097694767372
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 4.595688
$ws.Cells.Item(17, 14).Value = 13.787064
$ws.Cells.Item(17, 15).Value = 0.05841663420406906
$ws.Cells.Item(17, 16).Value = 0.05841663420406905
$ws.Cells.Item(17, 17).Value = 3.604817837904001
$ws.Cells.Item(17, 18).Value = 32.443360541136
$ws.Cells.Item(17, 19).Value = 0.01017674837829775
$ws.Cells.Item(17, 20).Value = 0.01017674837829775
